$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "contents" sheet: add a new row describing the "demo2" tab
# ------------------------------------------------------------------
$contents = $wb.Worksheets.Item("contents")
$contents.Cells.Item(3, 1).Value = "TestTab"
$contents.Cells.Item(3, 2).Value = "Display"
$contents.Cells.Item(3, 3).Value = "demo2"
$contents.Cells.Item(3, 4).Value = "users"

# ------------------------------------------------------------------
# 2. "demographics" sheet: add the missing "type" header in A1
# ------------------------------------------------------------------
$demographics = $wb.Worksheets.Item("demographics")
$demographics.Cells.Item(1, 1).Value = "type"

# ------------------------------------------------------------------
# 3. add a brand new "demo2" sheet, after "demographics", and fill
#    it with the same kind of "display" data
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$demo2 = $wb.Worksheets.Add($null, $lastSheet)
$demo2.Name = "demo2"

$headers = @("type", "name", "value", "parameter_list", "variable", "variable_value")
for ($c = 1; $c -le 6; $c++) {
    $demo2.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$genderText = 'text = "Gender", colour = "blue"'
$urnText = 'text = "Urn", colour = "blue"'

$rows = @(
    @("box", "box1", "table, plot", $genderText, "gender"),
    @("box", "box2", "table, plot", $genderText, "gender"),
    @("box", "box3", "table, plot", $genderText, "gender"),
    @("box", "box4", "table, plot", $genderText, "gender"),
    @("box", "box5", "table, plot", $genderText, "gender"),
    @("box", "box6", "table, plot", $urnText, "avatar"),
    @("box", "box7", "table, plot", $urnText, "avatar"),
    @("box", "box8", "table, plot", $urnText, "avatar"),
    @("box", "box9", "table, plot", $urnText, "avatar")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($c = 1; $c -le 5; $c++) {
        $demo2.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

# match the bold/centred header formatting used on the other sheets
# (copy the style instead of setting font/alignment individually so we
# don't create extra cellXfs entries)
$demographics.Range("B1").Copy()
$demo2.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 4. selections / active sheet bookkeeping
# ------------------------------------------------------------------
$demo2.Range("E2:E6").Select() | Out-Null

$demographics.Range("A2").Select() | Out-Null

$contents.Activate() | Out-Null
$contents.Range("D3").Select() | Out-Null
